$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.991.93"
$ws.Range("E2").Value = "  +1.44%  "

$ws.Range("D3").Value = "2.542.04"
$ws.Range("E3").Value = "  +2.38%  "

$c = $ws.Range("D4")
$c.Value = "'0.999"
$c.ClearFormats()
$ws.Range("E4").Value = "  -0.16%  "

$c = $ws.Range("D5")
$c.Value = "'526.32"
$c.ClearFormats()
$ws.Range("E5").Value = "  +1.38%  "

$c = $ws.Range("D6")
$c.Value = "'134.00"
$c.ClearFormats()
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("E7").Value = "  +0.32%  "

$ws.Range("E8").Value = "  +1.36%  "

$ws.Range("D9").Value = "2.540.07"
$ws.Range("E9").Value = "  +2.30%  "

$c = $ws.Range("D10")
$c.Value = "'0.0988"
$c.ClearFormats()
$ws.Range("E10").Value = "  +0.70%  "

$ws.Range("E11").Value = "  -1.49%  "

$c = $ws.Range("D12")
$c.Value = "'5.17"
$c.ClearFormats()
$ws.Range("E12").Value = "  -2.21%  "

$ws.Range("E13").Value = "  -0.51%  "

$ws.Range("D14").Value = "2.991.27"
$ws.Range("E14").Value = "  +2.80%  "

$ws.Range("D15").Value = "58.916.31"
$ws.Range("E15").Value = "  +1.38%  "

$c = $ws.Range("D16")
$c.Value = "'22.38"
$c.ClearFormats()
$ws.Range("E16").Value = "  +1.93%  "

$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("D18").Value = "2.542.36"
$ws.Range("E18").Value = "  +2.67%  "

$c = $ws.Range("D19")
$c.Value = "'10.73"
$c.ClearFormats()
$ws.Range("E19").Value = "  +1.13%  "

$c = $ws.Range("D20")
$c.Value = "'323.75"
$c.ClearFormats()
$ws.Range("E20").Value = "  +1.34%  "

$c = $ws.Range("D21")
$c.Value = "'4.20"
$c.ClearFormats()
$ws.Range("E21").Value = "  +0.72%  "

$c = $ws.Range("D22")
$c.Value = "'6.13"
$c.ClearFormats()
$ws.Range("E22").Value = "  +6.91%  "

$ws.Range("E23").Value = "  +0.03%  "

$c = $ws.Range("D24")
$c.Value = "'65.17"
$c.ClearFormats()
$ws.Range("E24").Value = "  +0.61%  "

$ws.Range("E25").Value = "  +0.12%  "

$c = $ws.Range("D26")
$c.Value = "'1.00"
$c.ClearFormats()
$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("E27").Value = "  -0.49%  "

$c = $ws.Range("D28")
$c.Value = "'7.43"
$c.ClearFormats()
$ws.Range("E28").Value = "  +1.21%  "

$ws.Range("D29").Value = "0.0₃0755"
$ws.Range("E29").Value = "  +1.09%  "

$ws.Range("E30").Value = "  +2.41%  "

$ws.Range("E31").Value = "  +2.46%  "

$c = $ws.Range("D32")
$c.Value = "'168.72"
$c.ClearFormats()
$ws.Range("E32").Value = "  -0.95%  "

$c = $ws.Range("D33")
$c.Value = "'6.36"
$c.ClearFormats()
$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("E34").Value = "  +0.00%  "

$c = $ws.Range("D35")
$c.Value = "'0.999"
$c.ClearFormats()
$ws.Range("E35").Value = "  +0.09%  "

$c = $ws.Range("D36")
$c.Value = "'18.31"
$c.ClearFormats()
$ws.Range("E36").Value = "  +1.57%  "

$ws.Range("E37").Value = "  -3.44%  "

$ws.Range("E38").Value = "  -0.30%  "

$ws.Range("E39").Value = "  +2.44%  "

$c = $ws.Range("D40")
$c.Value = "'36.70"
$c.ClearFormats()
$ws.Range("E40").Value = "  +0.19%  "

$ws.Range("E41").Value = "  -1.23%  "

$c = $ws.Range("D42")
$c.Value = "'280.33"
$c.ClearFormats()
$ws.Range("E42").Value = "  +2.13%  "

$ws.Range("E43").Value = "  +1.47%  "

$c = $ws.Range("D44")
$c.Value = "'5.11"
$c.ClearFormats()
$ws.Range("E44").Value = "  +1.43%  "

$c = $ws.Range("D45")
$c.Value = "'131.97"
$c.ClearFormats()
$ws.Range("E45").Value = "  +6.44%  "

$c = $ws.Range("D46")
$c.Value = "'0.604"
$c.ClearFormats()
$ws.Range("E46").Value = "  +1.78%  "

$c = $ws.Range("D47")
$c.Value = "'0.0921"
$c.ClearFormats()
$ws.Range("E47").Value = "  +1.45%  "

$ws.Range("E48").Value = "  +3.39%  "

$c = $ws.Range("D49")
$c.Value = "'17.87"
$c.ClearFormats()
$ws.Range("E49").Value = "  +1.40%  "

$c = $ws.Range("D50")
$c.Value = "'0.0216"
$c.ClearFormats()
$ws.Range("E50").Value = "  +1.55%  "

$c = $ws.Range("D51")
$c.Value = "'17.16"
$c.ClearFormats()
$ws.Range("E51").Value = "  +1.09%  "
